$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = 1
}

for ($row = 18; $row -le 21; $row++) {
    $ws.Cells.Item($row, 1).Value = 3
}

for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 6).Value = "no_pic"
}
